$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 136.5
$ws.Cells.Item(9, 9).Value = 54.75
$ws.Cells.Item(9, 11).Value = 54.75
$ws.Cells.Item(9, 13).Value = 114.25
$ws.Cells.Item(101, 8).Value = 313
$ws.Cells.Item(101, 9).Value = 291.4
$ws.Cells.Item(101, 10).Value = 367
$ws.Cells.Item(101, 11).Value = 874.1999999999999
$ws.Cells.Item(101, 12).Value = 1101
$ws.Cells.Item(101, 13).Value = 747.8000000000001
$ws.Cells.Item(101, 14).Value = -4345
$ws.Cells.Item(113, 8).Value = 3835.2856
$ws.Cells.Item(113, 9).Value = 3835.2856
$ws.Cells.Item(113, 11).Value = 3835.2856
$ws.Cells.Item(113, 13).Value = -581.2856000000002
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 9367.666999999999
$ws.Cells.Item(132, 9).Value = 12162.889
$ws.Cells.Item(132, 11).Value = 36488.667
$ws.Cells.Item(132, 13).Value = -33958.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 717.2857
$ws.Cells.Item(32, 9).Value = 717.2857
$ws.Cells.Item(32, 11).Value = 717.2857
$ws.Cells.Item(32, 13).Value = -430.2857
$ws.Cells.Item(45, 8).Value = 2320.9285
$ws.Cells.Item(45, 9).Value = 1304.875
$ws.Cells.Item(45, 11).Value = 1304.875
$ws.Cells.Item(45, 13).Value = -927.875
$ws.Cells.Item(74, 8).Value = 1548
$ws.Cells.Item(74, 9).Value = 1096
$ws.Cells.Item(74, 11).Value = 1096
$ws.Cells.Item(74, 13).Value = -222
$ws.Cells.Item(77, 8).Value = 1548
$ws.Cells.Item(77, 9).Value = 1096
$ws.Cells.Item(77, 11).Value = 5480
$ws.Cells.Item(77, 13).Value = -1112
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(53, 8).Value = 105389.5
$ws.Cells.Item(53, 10).Value = 105389.5
$ws.Cells.Item(53, 12).Value = 105389.5
$ws.Cells.Item(53, 14).Value = -106537.5
$ws.Cells.Item(86, 8).Value = 3533.348
$ws.Cells.Item(86, 9).Value = 1891.2
$ws.Cells.Item(86, 10).Value = 6612.375
$ws.Cells.Item(86, 11).Value = 1891.2
$ws.Cells.Item(86, 12).Value = 6612.375
$ws.Cells.Item(86, 13).Value = -768.2
$ws.Cells.Item(86, 14).Value = -8858.375
$ws.Cells.Item(89, 8).Value = 3533.348
$ws.Cells.Item(89, 9).Value = 1891.2
$ws.Cells.Item(89, 10).Value = 6612.375
$ws.Cells.Item(89, 11).Value = 9456
$ws.Cells.Item(89, 12).Value = 33061.875
$ws.Cells.Item(89, 13).Value = -3840
$ws.Cells.Item(89, 14).Value = -44293.875
$ws.Cells.Item(107, 8).Value = 5298.5
$ws.Cells.Item(107, 9).Value = 3275.25
$ws.Cells.Item(107, 10).Value = 7996.1665
$ws.Cells.Item(107, 11).Value = 3275.25
$ws.Cells.Item(107, 12).Value = 7996.1665
$ws.Cells.Item(107, 13).Value = -1355.25
$ws.Cells.Item(107, 14).Value = -11836.1665
$ws.Cells.Item(116, 8).Value = 64580
$ws.Cells.Item(116, 10).Value = 64580
$ws.Cells.Item(116, 12).Value = 64580
$ws.Cells.Item(116, 14).Value = -73758
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1096.2
$ws.Cells.Item(16, 9).Value = 1138.7858
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 11).Value = 1138.7858
$ws.Cells.Item(16, 12).Value = 500
$ws.Cells.Item(16, 13).Value = -851.7858000000001
$ws.Cells.Item(16, 14).Value = -1074
$ws.Cells.Item(19, 8).Value = 130
$ws.Cells.Item(19, 9).Value = 130
$ws.Cells.Item(19, 11).Value = 130
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(24, 8).Value = 130
$ws.Cells.Item(24, 9).Value = 130
$ws.Cells.Item(24, 11).Value = 130
$ws.Cells.Item(24, 13).Value = 40
$ws.Cells.Item(31, 8).Value = 5010.8604
$ws.Cells.Item(31, 9).Value = 2778.5862
$ws.Cells.Item(31, 10).Value = 9634.857
$ws.Cells.Item(31, 11).Value = 2778.5862
$ws.Cells.Item(31, 12).Value = 9634.857
$ws.Cells.Item(31, 13).Value = -2483.5862
$ws.Cells.Item(31, 14).Value = -10224.857
$ws.Cells.Item(34, 8).Value = 5010.8604
$ws.Cells.Item(34, 9).Value = 2778.5862
$ws.Cells.Item(34, 10).Value = 9634.857
$ws.Cells.Item(34, 11).Value = 2778.5862
$ws.Cells.Item(34, 12).Value = 9634.857
$ws.Cells.Item(34, 13).Value = -2576.5862
$ws.Cells.Item(34, 14).Value = -10038.857
$ws.Cells.Item(99, 8).Value = 2300
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(113, 8).Value = 1096.2
$ws.Cells.Item(113, 9).Value = 1138.7858
$ws.Cells.Item(113, 10).Value = 500
$ws.Cells.Item(113, 11).Value = 1138.7858
$ws.Cells.Item(113, 12).Value = 500
$ws.Cells.Item(113, 13).Value = 1031.2142
$ws.Cells.Item(113, 14).Value = -4840
$ws.Cells.Item(122, 8).Value = 1379.875
$ws.Cells.Item(122, 9).Value = 1077
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 3231
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -781
$ws.Cells.Item(122, 14).Value = -15400
$ws.Cells.Item(126, 8).Value = 2300
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 2666.3333
$ws.Cells.Item(132, 9).Value = 2666.3333
$ws.Cells.Item(132, 11).Value = 7998.999899999999
$ws.Cells.Item(132, 13).Value = -5468.999899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 286285.72
$ws.Cells.Item(4, 9).Value = 500500.25
$ws.Cells.Item(4, 10).Value = 666.3333
$ws.Cells.Item(4, 11).Value = 1501500.75
$ws.Cells.Item(4, 12).Value = 1998.9999
$ws.Cells.Item(4, 13).Value = -1501388.75
$ws.Cells.Item(4, 14).Value = -2222.9999
$ws.Cells.Item(6, 8).Value = 6
$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 18
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).ClearContents()
$ws.Cells.Item(6, 14).Value = 95
$ws.Cells.Item(12, 8).Value = 156
$ws.Cells.Item(12, 9).Value = 25.166666
$ws.Cells.Item(12, 11).Value = 75.49999800000001
$ws.Cells.Item(12, 13).Value = 97.50000199999999
$ws.Cells.Item(33, 8).Value = 415
$ws.Cells.Item(33, 9).Value = 398
$ws.Cells.Item(33, 10).Value = 500
$ws.Cells.Item(33, 11).Value = 2388
$ws.Cells.Item(33, 12).Value = 3000
$ws.Cells.Item(33, 13).Value = -2105
$ws.Cells.Item(33, 14).Value = -3566
$ws.Cells.Item(34, 8).Value = 585.2857
$ws.Cells.Item(34, 10).Value = 862.25
$ws.Cells.Item(34, 12).Value = 2586.75
$ws.Cells.Item(34, 14).Value = -2754.75
$ws.Cells.Item(38, 8).Value = 372.5
$ws.Cells.Item(38, 9).Value = 340.625
$ws.Cells.Item(38, 10).Value = 500
$ws.Cells.Item(38, 11).Value = 1021.875
$ws.Cells.Item(38, 12).Value = 1500
$ws.Cells.Item(38, 13).Value = -674.875
$ws.Cells.Item(38, 14).Value = -2194
$ws.Cells.Item(39, 8).Value = 6699
$ws.Cells.Item(39, 10).Value = 6699
$ws.Cells.Item(39, 12).Value = 20097
$ws.Cells.Item(39, 14).Value = -20685
$ws.Cells.Item(55, 8).Value = 6958.25
$ws.Cells.Item(55, 10).Value = 8777.666999999999
$ws.Cells.Item(55, 12).Value = 26333.001
$ws.Cells.Item(55, 14).Value = -26687.001
$ws.Cells.Item(68, 8).Value = 2000
$ws.Cells.Item(68, 9).Value = 1500
$ws.Cells.Item(68, 10).Value = 2500
$ws.Cells.Item(68, 11).Value = 4500
$ws.Cells.Item(68, 12).Value = 7500
$ws.Cells.Item(68, 13).Value = -3689
$ws.Cells.Item(68, 14).Value = -9122
$ws.Cells.Item(71, 8).Value = 2000
$ws.Cells.Item(71, 9).Value = 1500
$ws.Cells.Item(71, 10).Value = 2500
$ws.Cells.Item(71, 11).Value = 13500
$ws.Cells.Item(71, 12).Value = 22500
$ws.Cells.Item(71, 13).Value = -9444
$ws.Cells.Item(71, 14).Value = -30612
$ws.Cells.Item(114, 8).Value = 987.8
$ws.Cells.Item(114, 9).Value = 468
$ws.Cells.Item(114, 11).Value = 1404
$ws.Cells.Item(114, 13).Value = 1850
$ws.Cells.Item(129, 8).Value = 900
$ws.Cells.Item(129, 9).Value = 900
$ws.Cells.Item(129, 11).Value = 2700
$ws.Cells.Item(129, 13).Value = 2300
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1393.5
$ws.Cells.Item(97, 10).Value = 1999.5
$ws.Cells.Item(97, 12).Value = 1999.5
$ws.Cells.Item(97, 14).Value = -2991.5
$ws.Cells.Item(122, 8).Value = 3241.8572
$ws.Cells.Item(122, 9).Value = 2638.6
$ws.Cells.Item(122, 10).Value = 4750
$ws.Cells.Item(122, 11).Value = 7915.799999999999
$ws.Cells.Item(122, 12).Value = 14250
$ws.Cells.Item(122, 13).Value = -5465.799999999999
$ws.Cells.Item(122, 14).Value = -19150
$ws.Cells.Item(127, 8).Value = 15000
$ws.Cells.Item(127, 10).Value = 15000
$ws.Cells.Item(127, 12).Value = 15000
$ws.Cells.Item(127, 14).Value = -24920
$ws.Cells.Item(132, 8).Value = 4133
$ws.Cells.Item(132, 9).Value = 4133
$ws.Cells.Item(132, 11).Value = 12399
$ws.Cells.Item(132, 13).Value = -9869
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 6450
$ws.Cells.Item(61, 9).Value = 4900.25
$ws.Cells.Item(61, 11).Value = 4900.25
$ws.Cells.Item(61, 13).Value = -4698.25
$ws.Cells.Item(113, 8).Value = 6450
$ws.Cells.Item(113, 9).Value = 4900.25
$ws.Cells.Item(113, 11).Value = 4900.25
$ws.Cells.Item(113, 13).Value = -2730.25
$ws.Cells.Item(132, 8).Value = 5599.4
$ws.Cells.Item(132, 9).Value = 5332.3335
$ws.Cells.Item(132, 11).Value = 15997.0005
$ws.Cells.Item(132, 13).Value = -13467.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 48000
$ws.Cells.Item(16, 10).Value = 48000
$ws.Cells.Item(16, 12).Value = 48000
$ws.Cells.Item(16, 14).Value = -48584
$ws.Cells.Item(100, 8).Value = 948.2222
$ws.Cells.Item(100, 9).Value = 769.1429000000001
$ws.Cells.Item(100, 11).Value = 1538.2858
$ws.Cells.Item(100, 13).Value = -997.2858000000001
$ws.Cells.Item(107, 8).Value = 367.33334
$ws.Cells.Item(107, 9).Value = 322.14285
$ws.Cells.Item(107, 11).Value = 966.4285500000001
$ws.Cells.Item(107, 13).Value = 953.5714499999999
$ws.Cells.Item(113, 8).Value = 461.42856
$ws.Cells.Item(113, 9).Value = 324
$ws.Cells.Item(113, 11).Value = 972
$ws.Cells.Item(113, 13).Value = 1198
$ws.Cells.Item(122, 8).Value = 3073.5386
$ws.Cells.Item(122, 9).Value = 1996
$ws.Cells.Item(122, 11).Value = 5988
$ws.Cells.Item(122, 13).Value = -3538

Write-Output "applied 240 cell updates"